$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New FCF rows (46, 47) with their labels -----------------------------
# Add labels first so the new shared-string entries land before the
# "day 1 of project start" rename below (matches target sharedStrings order).
$ws.Range("A46").Value = "FCF project"
$ws.Range("A47").Value = "FCF owners"

# --- FCF project row formulas --------------------------------------------
$ws.Range("C46").Formula = '=C14-(C20-B20)-(C22-B22)+(C33-B33)+C11-C24'
$ws.Range("D46").Formula = '=D14-(D20-C20)-(D22-C22)+(D33-C33)+D11'
$ws.Range("E46").Formula = '=E14-(E20-D20)-(E22-D22)+(E33-D33)+E11'
$ws.Range("F46").Formula = '=F14-(F20-E20)-(F22-E22)+(F33-E33)+F11'
$ws.Range("G46").Formula = '=G14-(G20-F20)-(G22-F22)+(G33-F33)+G11'
$ws.Range("H46").Formula = '=H14-(H20-G20)-(H22-G22)+(H33-G33)+H11'

# --- FCF owners row formulas ----------------------------------------------
$ws.Range("C47").Formula = '=-(C27-B27)+C14+C9'
$ws.Range("D47").Formula = '=-(D27-C27)-D41'
$ws.Range("E47").Formula = '=-(E27-D27)-E41'
$ws.Range("F47").Formula = '=-(F27-E27)-F41'
$ws.Range("G47").Formula = '=-(G27-F27)-G41'
$ws.Range("H47").Formula = '=-(H27-G27)-H41'

# --- Carry the same "helper" number style (s=16) onto the blank cells ----
# that bound the new block, mirroring the formatting left behind by the
# original author's row fill/selection.
$ws.Range("B35").Copy()
$ws.Range("B46").PasteSpecial(-4122)
$ws.Range("B47").PasteSpecial(-4122)
$ws.Range("I47").PasteSpecial(-4122)
$ws.Range("E50").PasteSpecial(-4122)

# --- Move "minimal cash balance" amount from B43 to C43 -------------------
$ws.Range("A43").Copy()
$ws.Range("C43").PasteSpecial(-4122)
$ws.Range("B43").Clear()
$ws.Range("C43").Value = 500

# --- Point the "decrease of assets on bank accounts" row at the moved ----
# minimal-cash-balance cell (was $B$43, now $C$43).
$ws.Range("C41").Formula = '=IF(C38>=0,MIN(C38,B24-$C$43),(C38+C39))'
$ws.Range("D41").Formula = '=IF(D38>=0,MIN(D38,C24-$C$43),(D38+D39))'
$ws.Range("E41").Formula = '=IF(E38>=0,MIN(E38,D24-$C$43),(E38+E39))'
$ws.Range("F41").Formula = '=IF(F38>=0,MIN(F38,E24-$C$43),(F38+F39))'
$ws.Range("G41").Formula = '=IF(G38>=0,MIN(G38,F24-$C$43),(G38+G39))'
$ws.Range("H41").Formula = '=IF(H38>=0,MIN(H38,G24-$C$43),(H38+H39))'

# --- Core model edits: the project now draws down to 0 at day 1 ----------
$ws.Range("B24").Value = 0
$ws.Range("B27").Formula = '=B24'

# --- Rename the "day of project start" label ------------------------------
$ws.Range("B18").Value = "day 1 of project start"

# --- Restore the view the author left the sheet in ------------------------
$ws.Range("D47").Select()
$window = $excel.ActiveWindow
$window.ScrollRow = 12
